# Update InsideBet Data: Automatizado
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Arsenal
$ws.Range("C2").Value = 27
$ws.Range("E2").Value = 7
$ws.Range("G2").Value = 52
$ws.Range("H2").Value = 20
$ws.Range("J2").Value = 58
$ws.Range("K2").Value = 2.15
$ws.Range("L2").Value = "L W W D D"

# Row 21 - Wolves
$ws.Range("C21").Value = 27
$ws.Range("E21").Value = 7
$ws.Range("G21").Value = 18
$ws.Range("H21").Value = 50
$ws.Range("J21").Value = 10
$ws.Range("K21").Value = 0.37
$ws.Range("L21").Value = "L L L D D"
$ws.Range("M21").Value = 29798
